$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'323.08"
$ws.Range("E2").Value = "'3.41%"
$ws.Range("D3").Value = "'39.83"
$ws.Range("E3").Value = "'6.43%"
$ws.Range("D4").Value = "'5.230"
$ws.Range("E4").Value = "'1.89%"
$ws.Range("D5").Value = "'0.08105"
$ws.Range("E5").Value = "'2.53%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'8.621"
$ws.Range("E6").Value = "'4.29%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.918"
$ws.Range("E7").Value = "'0.42%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.958"
$ws.Range("E8").Value = "'-1.42%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9345"
$ws.Range("E9").Value = "'1.20%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1304"
$ws.Range("E10").Value = "'14.80%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1958"
$ws.Range("E11").Value = "'3.22%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09112"
$ws.Range("E12").Value = "'0.55%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03423"
$ws.Range("E13").Value = "'2.96%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09532"
$ws.Range("E14").Value = "'-0.78%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001402"
$ws.Range("E15").Value = "'1.59%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04437"
$ws.Range("E16").Value = "'1.70%"
$ws.Range("D17").Value = "'0.006439"
$ws.Range("E17").Value = "'4.07%"
$ws.Range("D18").Value = "'3.358"
$ws.Range("E18").Value = "'-5.94%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'4.528"
$ws.Range("E19").Value = "'2.53%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3535"
$ws.Range("E20").Value = "'2.52%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'6.716"
$ws.Range("E21").Value = "'14.00%"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1328"
$ws.Range("E22").Value = "'3.12%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.2311"
$ws.Range("E23").Value = "'-10.80%"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'-1.20%"
$ws.Range("D25").Value = "'0.004357"
$ws.Range("E25").Value = "'-6.18%"
$ws.Range("E26").Value = "'-5.26%"
$ws.Range("D27").Value = "'0.0003990"
$ws.Range("E27").Value = "'-0.08%"
$ws.Range("D39").Value = "'0.02461"
$ws.Range("E39").Value = "'8.25%"
$ws.Range("D40").Value = "'0.05235"
$ws.Range("E40").Value = "'2.89%"
$ws.Range("D41").Value = "'0.007644"
$ws.Range("E41").Value = "'2.40%"
$ws.Range("D42").Value = "'0.1432"
$ws.Range("E42").Value = "'5.65%"
$ws.Range("D43").Value = "'0.008726"
$ws.Range("E43").Value = "'-3.41%"
$ws.Range("E44").Value = "'8.10%"
$ws.Range("D45").Value = "'0.008130"
$ws.Range("E45").Value = "'-5.87%"
$ws.Range("D46").Value = "'0.00006596"
$ws.Range("E46").Value = "'-1.33%"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("D48").Value = "'0.002853"
$ws.Range("E48").Value = "'-13.32%"
$ws.Range("D49").Value = "'0.002483"
$ws.Range("E49").Value = "'148.10%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.20%"
